# support slti, sltui instructions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update ALUOp column (K) for rows that used the "branch" ALUOp code:
#     it moves from 2'bXXX to 3'bXXX now that a 3rd ALUOp bit is needed
#     to distinguish the new slti/sltui control signal.
$aluOpRows = @(2,3,4,5,6,7,8,9,10,11)
foreach ($r in $aluOpRows) {
    $ws.Cells.Item($r, 11).Value = "3'b010"
}

$ws.Cells.Item(12, 11).Value = "3'bxxx"

$ws.Cells.Item(15, 11).Value = "3'b000"
$ws.Cells.Item(16, 11).Value = "3'b011"
$ws.Cells.Item(17, 11).Value = "3'b100"

# row 18 (lui) keeps its special row formatting, but the ALUOp cell itself
# switches back to the plain (non-bold) cell style, matching the rest of
# the ALUOp column.
$ws.Cells.Item(18, 11).Value = "3'bxxx"
$ws.Cells.Item(18, 11).Font.Bold = $false
$ws.Cells.Item(18, 11).HorizontalAlignment = -4108

$memOpRows = @(19,20,21,22,23,24,25,26)
foreach ($r in $memOpRows) {
    $ws.Cells.Item($r, 11).Value = "3'b000"
}

$ws.Cells.Item(27, 11).Value = "3'b001"
$ws.Cells.Item(28, 11).Value = "3'b001"

$ws.Cells.Item(29, 11).Value = "3'bxxx"
$ws.Cells.Item(30, 11).Value = "3'bxxx"

# --- Fill in the new slti (row 13) and sltui (row 14) instruction rows.
# slti
$ws.Cells.Item(13, 2).Value = "b'001010"
$ws.Cells.Item(13, 3).Value = "2'b00"
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = "2'b00"
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 11).Value = "3'b101"
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = "2'bxx"
$ws.Cells.Item(13, 14).Value = "x"

# sltui
$ws.Cells.Item(14, 2).Value = "b'001001"
$ws.Cells.Item(14, 3).Value = "2'b00"
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = "2'b00"
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 11).Value = "3'b101"
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = "2'bxx"
$ws.Cells.Item(14, 14).Value = "x"

# --- Sheet view: zoom out a bit and move the selection/scroll position.
$excel.ActiveWindow.Zoom = 145
$null = $ws.Range("K18").Select()
